# Apply the "monthly budget file for 13, 14, 15, 16, 17 lessons" edit.

$wb = $excel.ActiveWorkbook

# 1) Rename "Sheet2" to "Shortcut Keys"
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Name = "Shortcut Keys"

# 2) Fill in the new "Lesson Name" values for Episodes 13-17 on the
#    "Completed Lessons" sheet (column C, rows 14-18).
$ws = $wb.Worksheets.Item("Completed Lessons")

$ws.Range("C14").Value = "Entering Text to Create Spreadsheet Titles"
$ws.Range("C15").Value = "Working with Numeric Data in Excel"
$ws.Range("C16").Value = "Entering Date Values in Excel"
$ws.Range("C17").Value = "Working with Cell References"
$ws.Range("C18").Value = "Creating Basic Formulas in Excel"

# 3) Update the active selection on "Completed Lessons" to C19.
$ws.Activate()
$ws.Range("C19").Select()
